$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Cells that hold the SQL query text which needs its JOIN conditions
# updated to use the renamed id columns (study_id / participant_id).
$cells = @("B2", "C2", "B3", "B4", "B5", "B6", "B7")

foreach ($cellRef in $cells) {
    $rng = $ws.Range($cellRef)
    $text = $rng.Value2

    $text = $text.Replace(
        'df_participant prt ON std.id = prt."study.id"',
        'df_participant prt ON std.study_id = prt."study.study_id"'
    )
    $text = $text.Replace(
        'df_diagnoses dgn ON prt.id = dgn."participant.id"',
        'df_diagnoses dgn ON prt.participant_id = dgn."participant.participant_id"'
    )
    $text = $text.Replace(
        'df_treatments trt ON prt.id = trt."participant.id"',
        'df_treatments trt ON prt.participant_id = trt."participant.participant_id"'
    )
    $text = $text.Replace(
        'df_treatment_resp trr ON prt.id = trr."participant.id"',
        'df_treatment_resp trr ON prt.participant_id = trr."participant.participant_id"'
    )
    $text = $text.Replace(
        'df_survival srv ON prt.id = srv."participant.id"',
        'df_survival srv ON prt.participant_id = srv."participant.participant_id"'
    )
    $text = $text.Replace(
        'df_reference_files rfs ON std.id = rfs."study.id"',
        'df_reference_files rfs ON std.study_id = rfs."study.study_id"'
    )

    $rng.Value2 = $text
}

# Widen column C to fit the updated (longer) query text, matching the
# author's manual resize (bestFit cleared, explicit width set to 69).
$ws.Columns.Item(3).ColumnWidth = 68.17

# Update the sheet view/selection to match what was saved: no frozen
# top-left scroll position, selection parked on B2 instead of C7.
$ws.Range("B2").Select()
